# Applies the "Updated symbol list" data refresh to the crypto price sheet.
# Numeric-looking values (prices / percentages) are written with an explicit
# text NumberFormat first so Excel keeps them as literal text (matching the
# source data, which stores every cell as text) instead of silently parsing
# them into floating point numbers / percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" '275.26'
Set-TextValue "E2" '0.61%'
Set-TextValue "D3" '27.19'
Set-TextValue "E3" '2.09%'
Set-TextValue "E4" '-0.91%'
Set-TextValue "D5" '0.06394'
Set-TextValue "E5" '1.07%'
Set-TextValue "D6" '6.943'
Set-TextValue "E6" '0.53%'
Set-TextValue "D7" '1.218'
Set-TextValue "E7" '-2.66%'
Set-TextValue "D8" '0.8761'
Set-TextValue "E8" '-0.05%'
Set-TextValue "E9" '4.53%'
Set-TextValue "D10" '0.05073'
Set-TextValue "E10" '-1.09%'
Set-TextValue "D11" '0.07518'
Set-TextValue "E11" '2.89%'
Set-TextValue "D12" '0.02956'
Set-TextValue "E12" '-5.12%'
Set-TextValue "D13" '0.08995'
Set-TextValue "E13" '-0.42%'
Set-TextValue "D14" '0.001563'
Set-TextValue "E14" '-0.38%'
Set-TextValue "D15" '0.0006409'
Set-TextValue "E15" '1.19%'
Set-TextValue "D16" '0.006183'
Set-TextValue "E16" '2.82%'
Set-TextValue "D17" '3.466'
Set-TextValue "E17" '0.31%'
Set-TextValue "D18" '3.309'
Set-TextValue "E18" '-1.50%'
Set-TextValue "E19" '-0.52%'
Set-TextValue "E21" '1.04%'
Set-TextValue "D22" '3.904'
Set-TextValue "E22" '-0.10%'
Set-TextValue "D23" '0.04409'
Set-TextValue "E23" '-0.06%'
Set-TextValue "D25" '0.001175'
Set-TextValue "E25" '-0.21%'
Set-TextValue "D26" '0.003852'
Set-TextValue "E26" '-12.57%'
Set-TextValue "D27" '0.0001201'
Set-TextValue "E27" '0.02%'
Set-TextValue "E28" '14.06%'
Set-TextValue "D40" '0.04127'
Set-TextValue "E40" '2.68%'
Set-TextValue "D41" '0.006799'
Set-TextValue "E41" '2.18%'
Set-TextValue "D43" '0.002141'
Set-TextValue "E43" '2.41%'
Set-TextValue "E44" '-6.40%'
Set-TextValue "D45" '0.00005176'
Set-TextValue "E45" '-2.48%'
$ws.Range("B46").Value = 'BOLO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue "D46" '1.501'
Set-TextValue "E46" '-36.29%'
$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue "D47" '0.02001'
Set-TextValue "E47" '0.02%'
